# Auto-generated script to apply "Add data for 2025-09-30" update
# Updates the 2025 (column L) values across the Citywide Totals, By Neighborhood,
# and individual neighborhood sheets to reflect data through 2025-09-30.

$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5048
$ws.Range("L3").Value = 5433
$ws.Range("L5").Value = 4630
$ws.Range("L6").Value = 1334
$ws.Range("L7").Value = 323
$ws.Range("L8").Value = 12532
$ws.Range("L9").Value = 4561
$ws.Range("L10").Value = 43672
$ws.Range("L11").Value = 77816

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 687
$ws.Range("L4").Value = 460
$ws.Range("L5").Value = 208
$ws.Range("L6").Value = 571
$ws.Range("L7").Value = 1816
$ws.Range("L8").Value = 3410
$ws.Range("L9").Value = 325
$ws.Range("L10").Value = 852
$ws.Range("L11").Value = 1303
$ws.Range("L12").Value = 309
$ws.Range("L13").Value = 209
$ws.Range("L14").Value = 556
$ws.Range("L15").Value = 589
$ws.Range("L16").Value = 570
$ws.Range("L17").Value = 106
$ws.Range("L18").Value = 491
$ws.Range("L19").Value = 1797
$ws.Range("L20").Value = 1310
$ws.Range("L21").Value = 193
$ws.Range("L22").Value = 288
$ws.Range("L23").Value = 962
$ws.Range("L24").Value = 374
$ws.Range("L25").Value = 443
$ws.Range("L27").Value = 986
$ws.Range("L29").Value = 2520
$ws.Range("L31").Value = 742
$ws.Range("L33").Value = 1823
$ws.Range("L34").Value = 530
$ws.Range("L36").Value = 1058
$ws.Range("L37").Value = 1735
$ws.Range("L38").Value = 109
$ws.Range("L39").Value = 132
$ws.Range("L41").Value = 261
$ws.Range("L42").Value = 1884
$ws.Range("L43").Value = 1091
$ws.Range("L44").Value = 738
$ws.Range("L46").Value = 216
$ws.Range("L47").Value = 692
$ws.Range("L48").Value = 2085
$ws.Range("L49").Value = 1302
$ws.Range("L50").Value = 794
$ws.Range("L51").Value = 1089
$ws.Range("L52").Value = 1039
$ws.Range("L53").Value = 1341
$ws.Range("L54").Value = 2873
$ws.Range("L55").Value = 776
$ws.Range("L56").Value = 438
$ws.Range("L59").Value = 218
$ws.Range("L60").Value = 546
$ws.Range("L63").Value = 344
$ws.Range("L64").Value = 699
$ws.Range("L65").Value = 1054
$ws.Range("L67").Value = 1404
$ws.Range("L69").Value = 375
$ws.Range("L70").Value = 514
$ws.Range("L71").Value = 272
$ws.Range("L72").Value = 518
$ws.Range("L73").Value = 987
$ws.Range("L74").Value = 235
$ws.Range("L76").Value = 2394
$ws.Range("L78").Value = 1205
$ws.Range("L79").Value = 1560
$ws.Range("L82").Value = 176
$ws.Range("L83").Value = 1166
$ws.Range("L84").Value = 626
$ws.Range("L85").Value = 2579
$ws.Range("L86").Value = 696
$ws.Range("L87").Value = 321
$ws.Range("L88").Value = 649
$ws.Range("L89").Value = 1564
$ws.Range("L90").Value = 799
$ws.Range("L91").Value = 732
$ws.Range("L92").Value = 242
$ws.Range("L93").Value = 607
$ws.Range("L94").Value = 2006
$ws.Range("L95").Value = 911
$ws.Range("L96").Value = 993
$ws.Range("L97").Value = 1315
$ws.Range("L98").Value = 813
$ws.Range("L99").Value = 1112
$ws.Range("L100").Value = 215
$ws.Range("L101").Value = 77816

# Sheet: Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L5").Value = 49
$ws.Range("L8").Value = 95
$ws.Range("L11").Value = 556

# Sheet: West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L5").Value = 71
$ws.Range("L8").Value = 183
$ws.Range("L10").Value = 548
$ws.Range("L11").Value = 993

# Sheet: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 184
$ws.Range("L3").Value = 182
$ws.Range("L8").Value = 421
$ws.Range("L10").Value = 704
$ws.Range("L11").Value = 1816

# Sheet: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L10").Value = 785
$ws.Range("L11").Value = 1303

# Sheet: O'Hare
$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L10").Value = 418
$ws.Range("L11").Value = 514

# Sheet: Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 72
$ws.Range("L6").Value = 39
$ws.Range("L10").Value = 1077
$ws.Range("L11").Value = 1564

# Sheet: South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 254
$ws.Range("L3").Value = 348
$ws.Range("L5").Value = 227
$ws.Range("L7").Value = 20
$ws.Range("L8").Value = 471
$ws.Range("L10").Value = 1020
$ws.Range("L11").Value = 2579

# Sheet: Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L8").Value = 124
$ws.Range("L9").Value = 91
$ws.Range("L11").Value = 1039

# Sheet: Norwood Park
$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L10").Value = 268
$ws.Range("L11").Value = 375

# Sheet: Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L10").Value = 901
$ws.Range("L11").Value = 1341

# Sheet: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 378
$ws.Range("L8").Value = 551
$ws.Range("L9").Value = 284
$ws.Range("L10").Value = 1563
$ws.Range("L11").Value = 3410

# Sheet: Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L8").Value = 40
$ws.Range("L10").Value = 131
$ws.Range("L11").Value = 216

# Sheet: Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 38
$ws.Range("L5").Value = 33
$ws.Range("L10").Value = 334
$ws.Range("L11").Value = 546

# Sheet: Oakland
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L8").Value = 93
$ws.Range("L11").Value = 272

# Sheet: South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 115
$ws.Range("L3").Value = 150
$ws.Range("L5").Value = 94
$ws.Range("L7").Value = 8
$ws.Range("L9").Value = 83
$ws.Range("L10").Value = 383
$ws.Range("L11").Value = 1166

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 211
$ws.Range("L6").Value = 48
$ws.Range("L8").Value = 332
$ws.Range("L10").Value = 625
$ws.Range("L11").Value = 1823

# Sheet: Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L10").Value = 657
$ws.Range("L11").Value = 1560

# Sheet: Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L5").Value = 31
$ws.Range("L10").Value = 389
$ws.Range("L11").Value = 699

# Sheet: West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 72
$ws.Range("L6").Value = 14
$ws.Range("L10").Value = 345
$ws.Range("L11").Value = 911

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 192
$ws.Range("L8").Value = 315
$ws.Range("L9").Value = 176
$ws.Range("L10").Value = 655
$ws.Range("L11").Value = 1735

# Sheet: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L5").Value = 75
$ws.Range("L10").Value = 463
$ws.Range("L11").Value = 1054

# Sheet: Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L5").Value = 90
$ws.Range("L9").Value = 44
$ws.Range("L10").Value = 644
$ws.Range("L11").Value = 986

# Sheet: Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 83
$ws.Range("L3").Value = 123
$ws.Range("L10").Value = 536
$ws.Range("L11").Value = 1112

# Sheet: Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 66
$ws.Range("L10").Value = 420
$ws.Range("L11").Value = 742

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L8").Value = 228
$ws.Range("L11").Value = 1404

# Sheet: South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L10").Value = 294
$ws.Range("L11").Value = 626

# Sheet: West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 28
$ws.Range("L8").Value = 253
$ws.Range("L10").Value = 1473
$ws.Range("L11").Value = 2006

# Sheet: River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L8").Value = 145
$ws.Range("L9").Value = 119
$ws.Range("L10").Value = 1933
$ws.Range("L11").Value = 2394

# Sheet: Ukrainian Village
$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L8").Value = 43
$ws.Range("L10").Value = 213
$ws.Range("L11").Value = 321

# Sheet: East Side
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L9").Value = 16
$ws.Range("L10").Value = 164
$ws.Range("L11").Value = 443

# Sheet: Wrigleyville
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L2").Value = 4
$ws.Range("L10").Value = 158
$ws.Range("L11").Value = 215

# Sheet: Bucktown
$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L5").Value = 41
$ws.Range("L10").Value = 448
$ws.Range("L11").Value = 570

# Sheet: Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L10").Value = 1009
$ws.Range("L11").Value = 1302

# Sheet: West Town
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L8").Value = 228
$ws.Range("L10").Value = 855
$ws.Range("L11").Value = 1315

# Sheet: Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 52
$ws.Range("L6").Value = 13
$ws.Range("L10").Value = 417
$ws.Range("L11").Value = 776

# Sheet: Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L8").Value = 262
$ws.Range("L10").Value = 2185
$ws.Range("L11").Value = 2873

# Sheet: Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 49
$ws.Range("L10").Value = 651
$ws.Range("L11").Value = 987

# Sheet: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 276
$ws.Range("L3").Value = 342
$ws.Range("L5").Value = 165
$ws.Range("L9").Value = 234
$ws.Range("L10").Value = 946
$ws.Range("L11").Value = 2520

# Sheet: Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L8").Value = 180
$ws.Range("L10").Value = 1517
$ws.Range("L11").Value = 2085

# Sheet: Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 163
$ws.Range("L7").Value = 8
$ws.Range("L10").Value = 891
$ws.Range("L11").Value = 1797

# Sheet: Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L8").Value = 116
$ws.Range("L10").Value = 453
$ws.Range("L11").Value = 738

# Sheet: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 44
$ws.Range("L8").Value = 335
$ws.Range("L10").Value = 903
$ws.Range("L11").Value = 1884

# Sheet: Clearing
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L5").Value = 24
$ws.Range("L10").Value = 150
$ws.Range("L11").Value = 288

# Sheet: Ashburn
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L8").Value = 160
$ws.Range("L11").Value = 571

# Sheet: Hermosa
$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L8").Value = 45
$ws.Range("L10").Value = 121
$ws.Range("L11").Value = 261

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L5").Value = 54
$ws.Range("L10").Value = 467
$ws.Range("L11").Value = 1058

# Sheet: Boystown
$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L9").Value = 161
$ws.Range("L10").Value = 209

# Sheet: Avondale
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L8").Value = 92
$ws.Range("L10").Value = 611
$ws.Range("L11").Value = 852

# Sheet: Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L10").Value = 492
$ws.Range("L11").Value = 696

# Sheet: Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 69
$ws.Range("L10").Value = 752
$ws.Range("L11").Value = 1205

# Sheet: Dunning
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L10").Value = 222
$ws.Range("L11").Value = 374

# Sheet: Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L10").Value = 334
$ws.Range("L11").Value = 589

# Sheet: Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 46
$ws.Range("L8").Value = 239
$ws.Range("L10").Value = 510
$ws.Range("L11").Value = 962

# Sheet: Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L8").Value = 204
$ws.Range("L11").Value = 732

# Sheet: Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L8").Value = 222
$ws.Range("L10").Value = 611
$ws.Range("L11").Value = 1089

# Sheet: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L8").Value = 207
$ws.Range("L10").Value = 578
$ws.Range("L11").Value = 1310

# Sheet: Chinatown
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L2").Value = 7
$ws.Range("L11").Value = 193

# Sheet: Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 37
$ws.Range("L10").Value = 409
$ws.Range("L11").Value = 692

# Sheet: Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L8").Value = 206
$ws.Range("L9").Value = 46
$ws.Range("L10").Value = 355
$ws.Range("L11").Value = 799

# Sheet: Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L9").Value = 22
$ws.Range("L10").Value = 575
$ws.Range("L11").Value = 794

# Sheet: West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L8").Value = 119
$ws.Range("L11").Value = 607

# Sheet: Calumet Heights
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 43
$ws.Range("L3").Value = 43
$ws.Range("L5").Value = 28
$ws.Range("L11").Value = 491

# Sheet: Magnificent Mile
$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("L10").Value = 414
$ws.Range("L11").Value = 438

# Sheet: Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L9").Value = 25
$ws.Range("L11").Value = 208

# Sheet: Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L10").Value = 399
$ws.Range("L11").Value = 687

# Sheet: Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L8").Value = 51
$ws.Range("L11").Value = 518

# Sheet: Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L5").Value = 54
$ws.Range("L6").Value = 17
$ws.Range("L8").Value = 135
$ws.Range("L9").Value = 40
$ws.Range("L10").Value = 782
$ws.Range("L11").Value = 1091

# Sheet: Burnside
$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L5").Value = 7
$ws.Range("L6").Value = 3
$ws.Range("L11").Value = 106

# Sheet: Archer Heights
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L10").Value = 309
$ws.Range("L11").Value = 460

# Sheet: Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L8").Value = 107
$ws.Range("L10").Value = 291
$ws.Range("L11").Value = 530

# Sheet: Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L8").Value = 107
$ws.Range("L9").Value = 44
$ws.Range("L10").Value = 566
$ws.Range("L11").Value = 813

# Sheet: West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L8").Value = 51
$ws.Range("L11").Value = 242

# Sheet: Sheffield & DePaul
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("L10").Value = 114
$ws.Range("L11").Value = 176

# Sheet: Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L8").Value = 81
$ws.Range("L11").Value = 325

# Sheet: Montclare
$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L10").Value = 127
$ws.Range("L11").Value = 218

# Sheet: Greektown
$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("L8").Value = 105
$ws.Range("L9").Value = 132

# Sheet: United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 56
$ws.Range("L10").Value = 313
$ws.Range("L11").Value = 649

# Sheet: Printers Row
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("L10").Value = 202
$ws.Range("L11").Value = 235

# Sheet: Grant Park
$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("L8").Value = 8
$ws.Range("L10").Value = 109

# Sheet: Beverly
$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L8").Value = 51
$ws.Range("L10").Value = 193
$ws.Range("L11").Value = 309

